$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "3.00") are preserved exactly as text, not converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "39.610.46"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "2.161.17"
$ws.Range("E3").Value = "  +2.71%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "228.08"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("D7").Value = "63.38"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").Value = "0.0849"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D12").Value = "15.96"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "2.481.57"
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("D14").Value = "21.98"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "2.180.73"
$ws.Range("E17").Value = "  +4.21%  "
$ws.Range("D18").Value = "39.558.05"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "71.90"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "6.12"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "0.0₃0845"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "228.21"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("D26").Value = "9.63"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "172.10"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "19.77"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  +4.79%  "
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "4.69"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").Value = "6.97"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").Value = "0.0618"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").Value = "3.61"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "4.75"
$ws.Range("E40").Value = "  +13.04%  "
$ws.Range("D41").Value = "102.00"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "17.69"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "1.511.25"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "1.22"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").Value = "0.0925"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "7.73"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "3.00"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").Value = "0.000189"
$ws.Range("E51").Value = "  +34.37%  "
